# Arbeitszeit_Pichler.xlsx - add three new diary entries (rows 65-67) plus
# extend the "Theoretische Grundlagen" tag onto rows 55-58, then update the
# view selection to match the author's final cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- New shared strings & their first usages --------------------------------
# The order these new strings are introduced below reproduces the exact
# sharedStrings index allocation seen in the target workbook (67..70).
$ws.Range("I67").Value = "User Interface editiert und Download-Buttons hinzugefügt"
$ws.Range("I65").Value = "Download-Methoden revidiert bzw. Code inspiziert"
$ws.Range("I66").Value = "Weiterführende Arbeit am Benutzer-Interface, insbesondere an den Optionen"
$ws.Range("I55").Value = "Arbeit am Teil ""Theoretische Grundlagen"""
$ws.Range("I56").Value = "Arbeit am Teil ""Theoretische Grundlagen"""
$ws.Range("I57").Value = "Arbeit am Teil ""Theoretische Grundlagen"""
$ws.Range("I58").Value = "Arbeit am Teil ""Theoretische Grundlagen"""

# --- Three new diary rows (65-67) -------------------------------------------
$ws.Range("E65").Value = 43821
$ws.Range("F65").Value = 3
$ws.Range("G65").Value = "Stunden"
$ws.Range("H65").Value = "Programmieren"

$ws.Range("E66").Value = 43822
$ws.Range("F66").Value = 4
$ws.Range("G66").Value = "Stunden"
$ws.Range("H66").Value = "Programmieren"

$ws.Range("E67").Value = 43823
$ws.Range("F67").Value = 3
$ws.Range("G67").Value = "Stunden"
$ws.Range("H67").Value = "Programmieren"

# --- Recalculate so SUM(F:F) and the TODAY()-based formulas refresh ---------
$ws.Calculate()

# --- Update the view: scroll position & active selection --------------------
$win = $ws.Application.ActiveWindow
$win.ScrollRow = 52
$win.ScrollColumn = 5
$ws.Range("I60").Select()
